$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" by duplicating the current "2022-Q2"
#    sheet and placing the copy right before it (so order becomes:
#    总计, 2022-Q3, 2022-Q2, 2021-Q2).
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Update the figures on the new "2022-Q3" sheet (the numbers reported
#    for that quarter differ from the ones that used to live on the
#    "2022-Q2" sheet).
# ---------------------------------------------------------------------------
$q3.Cells.Item(2, 4).Value = "'20.44"
$q3.Cells.Item(2, 5).Value = "'90.19"
$q3.Cells.Item(2, 6).Value = "'5.94"
$q3.Cells.Item(2, 7).Value = "'1.2141"
$q3.Cells.Item(2, 8).Value = 8

$q3.Cells.Item(3, 4).Value = "'15.02"
$q3.Cells.Item(3, 5).Value = "'90.19"
$q3.Cells.Item(3, 6).Value = "'5.94"
$q3.Cells.Item(3, 7).Value = "'0.8922"
$q3.Cells.Item(3, 8).Value = 8

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: shift the existing quarters down one
#    row and insert the new 2022-Q3 figures on top.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Make room: copy row 3 down into row 4 (keeps formatting/style in sync)
$zj.Range("A3:D3").Copy($zj.Range("A4:D4"))

# Row 4 now holds what used to be 2021-Q2 (row 3 content, unchanged values)
$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 2).Value = "2021-Q2"
$zj.Cells.Item(4, 3).Value = 2
$zj.Cells.Item(4, 4).Value = 2.28

# Row 3 now becomes what used to be 2022-Q2 (row 2 content, unchanged values)
$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2022-Q2"
$zj.Cells.Item(3, 3).Value = 2
$zj.Cells.Item(3, 4).Value = 2.31

# Row 2 becomes the new 2022-Q3 entry
$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q3"
$zj.Cells.Item(2, 3).Value = 2
$zj.Cells.Item(2, 4).Value = 2.11

# ---------------------------------------------------------------------------
# 4. Keep the originally-active tab (the "2021-Q2" sheet) selected, since the
#    newly inserted "2022-Q3" sheet should not steal the active-tab marker.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
